$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
Set-TextValue 2 4 '34.350.44'
$ws.Cells.Item(2, 5).Value = '  +0.77%  '

# Row 3
Set-TextValue 3 4 '1.836.10'
$ws.Cells.Item(3, 5).Value = '  +3.32%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.06%  '

# Row 5
Set-TextValue 5 4 '224.96'
$ws.Cells.Item(5, 5).Value = '  -0.02%  '

# Row 6
Set-TextValue 6 4 '0.557'
$ws.Cells.Item(6, 5).Value = '  +1.59%  '

# Row 7
Set-TextValue 7 4 '0.999'
$ws.Cells.Item(7, 5).Value = '  -0.10%  '

# Row 8
Set-TextValue 8 4 '31.87'
$ws.Cells.Item(8, 5).Value = '  +1.08%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +4.55%  '

# Row 10
Set-TextValue 10 4 '0.0722'
$ws.Cells.Item(10, 5).Value = '  +10.24%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +0.49%  '

# Row 12
Set-TextValue 12 4 '2.105.10'
$ws.Cells.Item(12, 5).Value = '  +3.55%  '

# Row 13
Set-TextValue 13 4 '1.843.99'
$ws.Cells.Item(13, 5).Value = '  +3.72%  '

# Row 14
Set-TextValue 14 4 '0.648'
$ws.Cells.Item(14, 5).Value = '  +3.83%  '

# Row 15
Set-TextValue 15 4 '10.78'
$ws.Cells.Item(15, 5).Value = '  -2.26%  '

# Row 16
Set-TextValue 16 4 '34.390.36'
$ws.Cells.Item(16, 5).Value = '  +0.88%  '

# Row 17
Set-TextValue 17 4 '4.35'
$ws.Cells.Item(17, 5).Value = '  +3.67%  '

# Row 18
Set-TextValue 18 4 '69.85'
$ws.Cells.Item(18, 5).Value = '  +1.78%  '

# Row 19
Set-TextValue 19 4 '251.45'
$ws.Cells.Item(19, 5).Value = '  -0.91%  '

# Row 20
Set-TextValue 20 4 '0.0₃0797'

# Row 21
Set-TextValue 21 4 '11.31'
$ws.Cells.Item(21, 5).Value = '  +9.75%  '

# Row 22
Set-TextValue 22 4 '0.998'
$ws.Cells.Item(22, 5).Value = '  -0.14%  '

# Row 23
Set-TextValue 23 4 '4.28'
$ws.Cells.Item(23, 5).Value = '  +2.40%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +1.53%  '

# Row 25
Set-TextValue 25 4 '160.78'
$ws.Cells.Item(25, 5).Value = '  +2.87%  '

# Row 26
Set-TextValue 26 4 '16.67'
$ws.Cells.Item(26, 5).Value = '  +2.00%  '

# Row 27
Set-TextValue 27 4 '7.26'
$ws.Cells.Item(27, 5).Value = '  +4.16%  '

# Row 28
Set-TextValue 28 4 '0.115'
$ws.Cells.Item(28, 5).Value = '  +2.01%  '

# Row 30
Set-TextValue 30 4 '0.0536'
$ws.Cells.Item(30, 5).Value = '  +4.91%  '

# Row 31
Set-TextValue 31 4 '3.81'
$ws.Cells.Item(31, 5).Value = '  +1.44%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +1.79%  '

# Row 33
Set-TextValue 33 4 '3.60'
$ws.Cells.Item(33, 5).Value = '  +1.11%  '

# Row 34
Set-TextValue 34 4 '1.92'
$ws.Cells.Item(34, 5).Value = '  +4.48%  '

# Row 35
Set-TextValue 35 4 '1.452.89'
$ws.Cells.Item(35, 5).Value = '  +0.92%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +4.07%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +3.15%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +1.95%  '

# Row 39
Set-TextValue 39 4 '0.968'
$ws.Cells.Item(39, 5).Value = '  +9.57%  '

# Row 40
Set-TextValue 40 4 '82.10'
$ws.Cells.Item(40, 5).Value = '  -0.64%  '

# Row 41
Set-TextValue 41 4 '2.77'
$ws.Cells.Item(41, 5).Value = '  -2.59%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  +0.36%  '

# Row 43
Set-TextValue 43 4 '2.15'
$ws.Cells.Item(43, 5).Value = '  +5.23%  '

# Row 44
Set-TextValue 44 4 '6.08'
$ws.Cells.Item(44, 5).Value = '  +4.78%  '

# Row 45
Set-TextValue 45 4 '1.999.39'
$ws.Cells.Item(45, 5).Value = '  +3.32%  '

# Row 46
Set-TextValue 46 4 '0.0500'
$ws.Cells.Item(46, 5).Value = '  -2.35%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  +1.07%  '

# Row 48
Set-TextValue 48 4 '106.65'
$ws.Cells.Item(48, 5).Value = '  +8.65%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'PaxDollar'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 49 4 '0.999'
$ws.Cells.Item(49, 5).Value = '  -0.07%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 50 4 '11.96'
$ws.Cells.Item(50, 5).Value = '  -0.56%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +7.25%  '
